$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above the current row 177, pushing the existing
# rows 177-203 down to 180-206 (dimension grows from T203 to T206).
$ws.Range("A177:A179").EntireRow.Insert()

# New row 177 - Tuna, Especial quality, Provincia de Los Andes
$ws.Range("A177").Value = 3
$ws.Range("B177").Value = "Femacal de La Calera"
$ws.Range("C177").Value = "Coquimbo"
$ws.Range("D177").Value = 45015
$ws.Range("E177").Value = 5
$ws.Range("F177").Value = "Fruta"
$ws.Range("G177").Value = 100107
$ws.Range("H177").Value = "Otros"
$ws.Range("I177").Value = 100107011
$ws.Range("J177").Value = "Tuna"
$ws.Range("K177").Value = "Sin especificar"
$ws.Range("L177").Value = "Especial"
$ws.Range("M177").Value = 60
$ws.Range("N177").Value = 16000
$ws.Range("O177").Value = 16000
$ws.Range("P177").Value = 16000
$ws.Range("Q177").Value = "$/caja 16 kilos"
$ws.Range("R177").Value = "Provincia de Los Andes"
$ws.Range("S177").Value = 1000
$ws.Range("T177").Value = 16

# New row 178 - Tuna, Primera quality, Provincia de Los Andes
$ws.Range("A178").Value = 3
$ws.Range("B178").Value = "Femacal de La Calera"
$ws.Range("C178").Value = "Coquimbo"
$ws.Range("D178").Value = 45015
$ws.Range("E178").Value = 5
$ws.Range("F178").Value = "Fruta"
$ws.Range("G178").Value = 100107
$ws.Range("H178").Value = "Otros"
$ws.Range("I178").Value = 100107011
$ws.Range("J178").Value = "Tuna"
$ws.Range("K178").Value = "Sin especificar"
$ws.Range("L178").Value = "Primera"
$ws.Range("M178").Value = 67
$ws.Range("N178").Value = 14000
$ws.Range("O178").Value = 14000
$ws.Range("P178").Value = 14000
$ws.Range("Q178").Value = "$/caja 16 kilos"
$ws.Range("R178").Value = "Provincia de Los Andes"
$ws.Range("S178").Value = 875
$ws.Range("T178").Value = 16

# New row 179 - Tuna, Segunda quality, Provincia de Los Andes
$ws.Range("A179").Value = 3
$ws.Range("B179").Value = "Femacal de La Calera"
$ws.Range("C179").Value = "Coquimbo"
$ws.Range("D179").Value = 45015
$ws.Range("E179").Value = 5
$ws.Range("F179").Value = "Fruta"
$ws.Range("G179").Value = 100107
$ws.Range("H179").Value = "Otros"
$ws.Range("I179").Value = 100107011
$ws.Range("J179").Value = "Tuna"
$ws.Range("K179").Value = "Sin especificar"
$ws.Range("L179").Value = "Segunda"
$ws.Range("M179").Value = 50
$ws.Range("N179").Value = 12000
$ws.Range("O179").Value = 12000
$ws.Range("P179").Value = 12000
$ws.Range("Q179").Value = "$/caja 16 kilos"
$ws.Range("R179").Value = "Provincia de Los Andes"
$ws.Range("S179").Value = 750
$ws.Range("T179").Value = 16
